$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 6753.793  # H40: 6916.75 -> 6753.793
$ws.Cells.Item(40, 9).Value = 3126.1428  # I40: 3126 -> 3126.1428
$ws.Cells.Item(40, 10).Value = 7908.0454  # J40: 8180.3335 -> 7908.0454
$ws.Cells.Item(40, 11).Value = 3126.1428  # K40: 3126 -> 3126.1428
$ws.Cells.Item(40, 12).Value = 7908.0454  # L40: 8180.3335 -> 7908.0454
$ws.Cells.Item(40, 13).Value = -2951.1428  # M40: -2951 -> -2951.1428
$ws.Cells.Item(40, 14).Value = -8258.045399999999  # N40: -8530.333500000001 -> -8258.045399999999
$ws.Cells.Item(43, 8).Value = 5499.3335  # H43: 6249.25 -> 5499.3335
$ws.Cells.Item(43, 10).Value = 3999.5  # J43: 0 -> 3999.5
$ws.Cells.Item(43, 12).Value = 3999.5  # L43: 0 -> 3999.5
$ws.Cells.Item(43, 14).Value = -4137.5  # N43: None -> -4137.5
$ws.Cells.Item(108, 8).Value = 0  # H108: 99995 -> 0
$ws.Cells.Item(108, 10).Value = 0  # J108: 99995 -> 0
$ws.Cells.Item(108, 12).Value = 0  # L108: 99995 -> 0
$ws.Cells.Item(108, 14).ClearContents()  # N108: delete (was -107675)
$ws.Cells.Item(114, 8).Value = 90000  # H114: 94997.5 -> 90000
$ws.Cells.Item(114, 10).Value = 0  # J114: 99995 -> 0
$ws.Cells.Item(114, 12).Value = 0  # L114: 99995 -> 0
$ws.Cells.Item(114, 14).ClearContents()  # N114: delete (was -108673)
$ws.Cells.Item(116, 8).Value = 9995  # H116: 0 -> 9995
$ws.Cells.Item(116, 9).Value = 10000  # I116: 0 -> 10000
$ws.Cells.Item(116, 10).Value = 9990  # J116: 0 -> 9990
$ws.Cells.Item(116, 11).Value = 10000  # K116: 0 -> 10000
$ws.Cells.Item(116, 12).Value = 9990  # L116: 0 -> 9990
$ws.Cells.Item(116, 13).Value = -6558  # M116: None -> -6558
$ws.Cells.Item(116, 14).Value = -16874  # N116: None -> -16874
$ws.Cells.Item(120, 8).Value = 0  # H120: 99995 -> 0
$ws.Cells.Item(120, 10).Value = 0  # J120: 99995 -> 0
$ws.Cells.Item(120, 12).Value = 0  # L120: 99995 -> 0
$ws.Cells.Item(120, 14).ClearContents()  # N120: delete (was -109671)
$ws.Cells.Item(128, 8).Value = 0  # H128: 99995 -> 0
$ws.Cells.Item(128, 10).Value = 0  # J128: 99995 -> 0
$ws.Cells.Item(128, 12).Value = 0  # L128: 99995 -> 0
$ws.Cells.Item(128, 14).ClearContents()  # N128: delete (was -109955)
$ws.Cells.Item(133, 8).Value = 0  # H133: 99995 -> 0
$ws.Cells.Item(133, 10).Value = 0  # J133: 99995 -> 0
$ws.Cells.Item(133, 12).Value = 0  # L133: 99995 -> 0
$ws.Cells.Item(133, 14).ClearContents()  # N133: delete (was -110115)
$ws.Cells.Item(134, 8).Value = 0  # H134: 99995 -> 0
$ws.Cells.Item(134, 10).Value = 0  # J134: 99995 -> 0
$ws.Cells.Item(134, 12).Value = 0  # L134: 99995 -> 0
$ws.Cells.Item(134, 14).ClearContents()  # N134: delete (was -110135)

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(108, 8).Value = 75825  # H108: 75925 -> 75825
$ws.Cells.Item(108, 10).Value = 75825  # J108: 75925 -> 75825
$ws.Cells.Item(108, 12).Value = 75825  # L108: 75925 -> 75825
$ws.Cells.Item(108, 14).Value = -83505  # N108: -83605 -> -83505
$ws.Cells.Item(115, 8).Value = 49999  # H115: 0 -> 49999
$ws.Cells.Item(115, 10).Value = 49999  # J115: 0 -> 49999
$ws.Cells.Item(115, 12).Value = 49999  # L115: 0 -> 49999
$ws.Cells.Item(115, 14).Value = -53133  # N115: None -> -53133
$ws.Cells.Item(118, 8).Value = 29999.5  # H118: 0 -> 29999.5
$ws.Cells.Item(118, 10).Value = 29999.5  # J118: 0 -> 29999.5
$ws.Cells.Item(118, 12).Value = 29999.5  # L118: 0 -> 29999.5
$ws.Cells.Item(118, 14).Value = -33313.5  # N118: None -> -33313.5
$ws.Cells.Item(128, 8).Value = 89999  # H128: 99995 -> 89999
$ws.Cells.Item(128, 10).Value = 89999  # J128: 99995 -> 89999
$ws.Cells.Item(128, 12).Value = 89999  # L128: 99995 -> 89999
$ws.Cells.Item(128, 14).Value = -99959  # N128: -109955 -> -99959

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 0  # H81: 99995 -> 0
$ws.Cells.Item(81, 10).Value = 0  # J81: 99995 -> 0
$ws.Cells.Item(81, 12).Value = 0  # L81: 99995 -> 0
$ws.Cells.Item(81, 14).ClearContents()  # N81: delete (was -102117)
$ws.Cells.Item(84, 8).Value = 0  # H84: 99995 -> 0
$ws.Cells.Item(84, 10).Value = 0  # J84: 99995 -> 0
$ws.Cells.Item(84, 12).Value = 0  # L84: 299985 -> 0
$ws.Cells.Item(84, 14).ClearContents()  # N84: delete (was -310593)
$ws.Cells.Item(107, 8).Value = 6000  # H107: 3231.2 -> 6000
$ws.Cells.Item(107, 9).Value = 0  # I107: 1385.3334 -> 0
$ws.Cells.Item(107, 11).Value = 0  # K107: 1385.3334 -> 0
$ws.Cells.Item(107, 13).ClearContents()  # M107: delete (was 534.6666)
$ws.Cells.Item(108, 8).Value = 0  # H108: 99995 -> 0
$ws.Cells.Item(108, 10).Value = 0  # J108: 99995 -> 0
$ws.Cells.Item(108, 12).Value = 0  # L108: 99995 -> 0
$ws.Cells.Item(108, 14).ClearContents()  # N108: delete (was -107675)
$ws.Cells.Item(111, 8).Value = 99000  # H111: 98997 -> 99000
$ws.Cells.Item(111, 10).Value = 99000  # J111: 98997 -> 99000
$ws.Cells.Item(111, 12).Value = 99000  # L111: 98997 -> 99000
$ws.Cells.Item(111, 14).Value = -107180  # N111: -107177 -> -107180
$ws.Cells.Item(112, 8).Value = 0  # H112: 99995 -> 0
$ws.Cells.Item(112, 10).Value = 0  # J112: 99995 -> 0
$ws.Cells.Item(112, 12).Value = 0  # L112: 99995 -> 0
$ws.Cells.Item(112, 14).ClearContents()  # N112: delete (was -102949)
$ws.Cells.Item(116, 8).Value = 24950  # H116: 77700 -> 24950
$ws.Cells.Item(116, 10).Value = 24950  # J116: 77700 -> 24950
$ws.Cells.Item(116, 12).Value = 24950  # L116: 77700 -> 24950
$ws.Cells.Item(116, 14).Value = -34128  # N116: -86878 -> -34128
$ws.Cells.Item(120, 8).Value = 0  # H120: 99995 -> 0
$ws.Cells.Item(120, 10).Value = 0  # J120: 99995 -> 0
$ws.Cells.Item(120, 12).Value = 0  # L120: 99995 -> 0
$ws.Cells.Item(120, 14).ClearContents()  # N120: delete (was -109671)
$ws.Cells.Item(125, 8).Value = 0  # H125: 99995 -> 0
$ws.Cells.Item(125, 10).Value = 0  # J125: 99995 -> 0
$ws.Cells.Item(125, 12).Value = 0  # L125: 99995 -> 0
$ws.Cells.Item(125, 14).ClearContents()  # N125: delete (was -109835)
$ws.Cells.Item(127, 8).Value = 0  # H127: 99995 -> 0
$ws.Cells.Item(127, 10).Value = 0  # J127: 99995 -> 0
$ws.Cells.Item(127, 12).Value = 0  # L127: 99995 -> 0
$ws.Cells.Item(127, 14).ClearContents()  # N127: delete (was -109915)
$ws.Cells.Item(139, 8).Value = 0  # H139: 99995 -> 0
$ws.Cells.Item(139, 10).Value = 0  # J139: 99995 -> 0
$ws.Cells.Item(139, 12).Value = 0  # L139: 99995 -> 0
$ws.Cells.Item(139, 14).ClearContents()  # N139: delete (was -110275)

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 319  # H22: 819.7778 -> 319
$ws.Cells.Item(22, 9).Value = 319  # I22: 884.5 -> 319
$ws.Cells.Item(22, 10).Value = 0  # J22: 302 -> 0
$ws.Cells.Item(22, 11).Value = 319  # K22: 884.5 -> 319
$ws.Cells.Item(22, 12).Value = 0  # L22: 302 -> 0
$ws.Cells.Item(22, 13).Value = 31  # M22: -534.5 -> 31
$ws.Cells.Item(22, 14).ClearContents()  # N22: delete (was -1002)
$ws.Cells.Item(53, 8).Value = 0  # H53: 60684 -> 0
$ws.Cells.Item(53, 10).Value = 0  # J53: 60684 -> 0
$ws.Cells.Item(53, 12).Value = 0  # L53: 60684 -> 0
$ws.Cells.Item(53, 14).ClearContents()  # N53: delete (was -61898)
$ws.Cells.Item(100, 8).Value = 0  # H100: 99995 -> 0
$ws.Cells.Item(100, 10).Value = 0  # J100: 99995 -> 0
$ws.Cells.Item(100, 12).Value = 0  # L100: 99995 -> 0
$ws.Cells.Item(100, 14).ClearContents()  # N100: delete (was -102159)
$ws.Cells.Item(110, 8).Value = 10000  # H110: 0 -> 10000
$ws.Cells.Item(110, 10).Value = 10000  # J110: 0 -> 10000
$ws.Cells.Item(110, 12).Value = 10000  # L110: 0 -> 10000
$ws.Cells.Item(110, 14).Value = -18180  # N110: None -> -18180
$ws.Cells.Item(111, 8).Value = 77251  # H111: 77351 -> 77251
$ws.Cells.Item(111, 10).Value = 77251  # J111: 77351 -> 77251
$ws.Cells.Item(111, 12).Value = 77251  # L111: 77351 -> 77251
$ws.Cells.Item(111, 14).Value = -85431  # N111: -85531 -> -85431
$ws.Cells.Item(116, 8).Value = 99983  # H116: 99989 -> 99983
$ws.Cells.Item(116, 10).Value = 99983  # J116: 99989 -> 99983
$ws.Cells.Item(116, 12).Value = 99983  # L116: 99989 -> 99983
$ws.Cells.Item(116, 14).Value = -109161  # N116: -109167 -> -109161
$ws.Cells.Item(118, 8).Value = 80000  # H118: 79997 -> 80000
$ws.Cells.Item(118, 10).Value = 80000  # J118: 79997 -> 80000
$ws.Cells.Item(118, 12).Value = 80000  # L118: 79997 -> 80000
$ws.Cells.Item(118, 14).Value = -83314  # N118: -83311 -> -83314
$ws.Cells.Item(130, 8).Value = 0  # H130: 99995 -> 0
$ws.Cells.Item(130, 10).Value = 0  # J130: 99995 -> 0
$ws.Cells.Item(130, 12).Value = 0  # L130: 99995 -> 0
$ws.Cells.Item(130, 14).ClearContents()  # N130: delete (was -110035)
$ws.Cells.Item(135, 8).Value = 0  # H135: 80000 -> 0
$ws.Cells.Item(135, 10).Value = 0  # J135: 80000 -> 0
$ws.Cells.Item(135, 12).Value = 0  # L135: 80000 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: delete (was -90140)
$ws.Cells.Item(141, 8).Value = 95600  # H141: 98000 -> 95600
$ws.Cells.Item(141, 10).Value = 95600  # J141: 98000 -> 95600
$ws.Cells.Item(141, 12).Value = 95600  # L141: 98000 -> 95600
$ws.Cells.Item(141, 14).Value = -105960  # N141: -108360 -> -105960

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 449.77777  # H5: 467.5 -> 449.77777
$ws.Cells.Item(5, 9).Value = 456  # I5: 477.14285 -> 456
$ws.Cells.Item(5, 11).Value = 1368  # K5: 1431.42855 -> 1368
$ws.Cells.Item(5, 13).Value = -1256  # M5: -1319.42855 -> -1256
$ws.Cells.Item(41, 10).Value = 0  # J41: 200 -> 0
$ws.Cells.Item(41, 12).Value = 0  # L41: 600 -> 0
$ws.Cells.Item(41, 14).ClearContents()  # N41: delete (was -1276)
$ws.Cells.Item(135, 8).Value = 449.77777  # H135: 467.5 -> 449.77777
$ws.Cells.Item(135, 9).Value = 456  # I135: 477.14285 -> 456
$ws.Cells.Item(135, 11).Value = 4104  # K135: 4294.28565 -> 4104
$ws.Cells.Item(135, 13).Value = -1569  # M135: -1759.28565 -> -1569

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 52751  # H80: 102003 -> 52751
$ws.Cells.Item(80, 9).Value = 2999  # I80: 0 -> 2999
$ws.Cells.Item(80, 10).Value = 69335  # J80: 102003 -> 69335
$ws.Cells.Item(80, 11).Value = 2999  # K80: 0 -> 2999
$ws.Cells.Item(80, 12).Value = 69335  # L80: 102003 -> 69335
$ws.Cells.Item(80, 13).Value = -2001  # M80: None -> -2001
$ws.Cells.Item(80, 14).Value = -71331  # N80: -103999 -> -71331
$ws.Cells.Item(83, 8).Value = 52751  # H83: 102003 -> 52751
$ws.Cells.Item(83, 9).Value = 2999  # I83: 0 -> 2999
$ws.Cells.Item(83, 10).Value = 69335  # J83: 102003 -> 69335
$ws.Cells.Item(83, 11).Value = 14995  # K83: 0 -> 14995
$ws.Cells.Item(83, 12).Value = 346675  # L83: 510015 -> 346675
$ws.Cells.Item(83, 13).Value = -10003  # M83: None -> -10003
$ws.Cells.Item(83, 14).Value = -356659  # N83: -519999 -> -356659
$ws.Cells.Item(110, 8).Value = 0  # H110: 99995 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 99995 -> 0
$ws.Cells.Item(110, 12).Value = 0  # L110: 99995 -> 0
$ws.Cells.Item(110, 14).ClearContents()  # N110: delete (was -108175)
$ws.Cells.Item(116, 8).Value = 0  # H116: 99995 -> 0
$ws.Cells.Item(116, 10).Value = 0  # J116: 99995 -> 0
$ws.Cells.Item(116, 12).Value = 0  # L116: 99995 -> 0
$ws.Cells.Item(116, 14).ClearContents()  # N116: delete (was -109173)
$ws.Cells.Item(119, 8).Value = 0  # H119: 99995 -> 0
$ws.Cells.Item(119, 10).Value = 0  # J119: 99995 -> 0
$ws.Cells.Item(119, 12).Value = 0  # L119: 99995 -> 0
$ws.Cells.Item(119, 14).ClearContents()  # N119: delete (was -109671)
$ws.Cells.Item(128, 8).Value = 0  # H128: 99995 -> 0
$ws.Cells.Item(128, 10).Value = 0  # J128: 99995 -> 0
$ws.Cells.Item(128, 12).Value = 0  # L128: 99995 -> 0
$ws.Cells.Item(128, 14).ClearContents()  # N128: delete (was -109955)
$ws.Cells.Item(130, 8).Value = 0  # H130: 99995 -> 0
$ws.Cells.Item(130, 10).Value = 0  # J130: 99995 -> 0
$ws.Cells.Item(130, 12).Value = 0  # L130: 99995 -> 0
$ws.Cells.Item(130, 14).ClearContents()  # N130: delete (was -110035)
$ws.Cells.Item(135, 8).Value = 0  # H135: 99995 -> 0
$ws.Cells.Item(135, 10).Value = 0  # J135: 99995 -> 0
$ws.Cells.Item(135, 12).Value = 0  # L135: 99995 -> 0
$ws.Cells.Item(135, 14).ClearContents()  # N135: delete (was -110135)

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2800  # H22: 2807.8462 -> 2800
$ws.Cells.Item(22, 9).Value = 2400  # I22: 2562.75 -> 2400
$ws.Cells.Item(22, 11).Value = 2400  # K22: 2562.75 -> 2400
$ws.Cells.Item(22, 13).Value = -2105  # M22: -2267.75 -> -2105
$ws.Cells.Item(27, 8).Value = 2800  # H27: 2807.8462 -> 2800
$ws.Cells.Item(27, 9).Value = 2400  # I27: 2562.75 -> 2400
$ws.Cells.Item(27, 11).Value = 2400  # K27: 2562.75 -> 2400
$ws.Cells.Item(27, 13).Value = -2293  # M27: -2455.75 -> -2293
$ws.Cells.Item(55, 8).Value = 2865.6667  # H55: 3186.125 -> 2865.6667
$ws.Cells.Item(55, 9).Value = 5500  # I55: 2999 -> 5500
$ws.Cells.Item(55, 10).Value = 2536.375  # J55: 3248.5 -> 2536.375
$ws.Cells.Item(55, 11).Value = 5500  # K55: 2999 -> 5500
$ws.Cells.Item(55, 12).Value = 2536.375  # L55: 3248.5 -> 2536.375
$ws.Cells.Item(55, 13).Value = -5327  # M55: -2826 -> -5327
$ws.Cells.Item(55, 14).Value = -2882.375  # N55: -3594.5 -> -2882.375
$ws.Cells.Item(82, 8).Value = 2252.6875  # H82: 2296.3333 -> 2252.6875
$ws.Cells.Item(82, 9).Value = 1985.909  # I82: 1986 -> 1985.909
$ws.Cells.Item(82, 10).Value = 2839.6  # J82: 3149.75 -> 2839.6
$ws.Cells.Item(82, 11).Value = 1985.909  # K82: 1986 -> 1985.909
$ws.Cells.Item(82, 12).Value = 2839.6  # L82: 3149.75 -> 2839.6
$ws.Cells.Item(82, 13).Value = -1624.909  # M82: -1625 -> -1624.909
$ws.Cells.Item(82, 14).Value = -3561.6  # N82: -3871.75 -> -3561.6
$ws.Cells.Item(85, 8).Value = 2252.6875  # H85: 2296.3333 -> 2252.6875
$ws.Cells.Item(85, 9).Value = 1985.909  # I85: 1986 -> 1985.909
$ws.Cells.Item(85, 10).Value = 2839.6  # J85: 3149.75 -> 2839.6
$ws.Cells.Item(85, 11).Value = 1985.909  # K85: 1986 -> 1985.909
$ws.Cells.Item(85, 12).Value = 2839.6  # L85: 3149.75 -> 2839.6
$ws.Cells.Item(85, 13).Value = -737.9090000000001  # M85: -738 -> -737.9090000000001
$ws.Cells.Item(85, 14).Value = -5335.6  # N85: -5645.75 -> -5335.6
$ws.Cells.Item(135, 8).Value = 99995  # H135: 0 -> 99995
$ws.Cells.Item(135, 10).Value = 99995  # J135: 0 -> 99995
$ws.Cells.Item(135, 12).Value = 99995  # L135: 0 -> 99995
$ws.Cells.Item(135, 14).Value = -110135  # N135: None -> -110135

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 99995  # H16: 0 -> 99995
$ws.Cells.Item(16, 10).Value = 99995  # J16: 0 -> 99995
$ws.Cells.Item(16, 12).Value = 99995  # L16: 0 -> 99995
$ws.Cells.Item(16, 14).Value = -100579  # N16: None -> -100579
$ws.Cells.Item(46, 8).Value = 0  # H46: 99995 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 99995 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 99995 -> 0
$ws.Cells.Item(46, 14).ClearContents()  # N46: delete (was -100457)
$ws.Cells.Item(109, 8).Value = 33333  # H109: 33014.332 -> 33333
$ws.Cells.Item(109, 10).Value = 33333  # J109: 33014.332 -> 33333
$ws.Cells.Item(109, 12).Value = 33333  # L109: 33014.332 -> 33333
$ws.Cells.Item(109, 14).Value = -36107  # N109: -35788.332 -> -36107
$ws.Cells.Item(130, 8).Value = 36666.668  # H130: 50000 -> 36666.668
$ws.Cells.Item(130, 9).Value = 35000  # I130: 50000 -> 35000
$ws.Cells.Item(130, 10).Value = 40000  # J130: 0 -> 40000
$ws.Cells.Item(130, 11).Value = 35000  # K130: 50000 -> 35000
$ws.Cells.Item(130, 12).Value = 40000  # L130: 0 -> 40000
$ws.Cells.Item(130, 13).Value = -29980  # M130: -44980 -> -29980
$ws.Cells.Item(130, 14).Value = -50040  # N130: None -> -50040
$ws.Cells.Item(134, 8).Value = 0  # H134: 99995 -> 0
$ws.Cells.Item(134, 10).Value = 0  # J134: 99995 -> 0
$ws.Cells.Item(134, 12).Value = 0  # L134: 299985 -> 0
$ws.Cells.Item(134, 14).ClearContents()  # N134: delete (was -305055)
